# Applies the numeric corrections from the scheduled-runner sheet update.
# Each worksheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) gets a batch of
# literal value writes; a couple of cells are fully cleared/added to match
# the row layout exactly (ClearContents for removed cells).
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1300.9697
$ws.Range("J17").Value = 1300.9697
$ws.Range("L17").Value = 3902.9091
$ws.Range("N17").Value = -4238.909100000001

$ws.Range("H40").Value = 4976.769
$ws.Range("I40").Value = 3596.25
$ws.Range("J40").Value = 7185.6
$ws.Range("K40").Value = 3596.25
$ws.Range("L40").Value = 7185.6
$ws.Range("M40").Value = -3421.25
$ws.Range("N40").Value = -7535.6

$ws.Range("H76").Value = 8490.357
$ws.Range("I76").Value = 7580.5
$ws.Range("J76").Value = 9172.75
$ws.Range("K76").Value = 7580.5
$ws.Range("L76").Value = 9172.75
$ws.Range("M76").Value = -7265.5
$ws.Range("N76").Value = -9802.75

$ws.Range("H79").Value = 8490.357
$ws.Range("I79").Value = 7580.5
$ws.Range("J79").Value = 9172.75
$ws.Range("K79").Value = 7580.5
$ws.Range("L79").Value = 9172.75
$ws.Range("M79").Value = -6488.5
$ws.Range("N79").Value = -11356.75

$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H86").Value = 2839.7334
$ws.Range("I86").Value = 2543.2727
$ws.Range("K86").Value = 2543.2727
$ws.Range("M86").Value = -1420.2727

$ws.Range("H89").Value = 2839.7334
$ws.Range("I89").Value = 2543.2727
$ws.Range("K89").Value = 12716.3635
$ws.Range("M89").Value = -7100.363499999999

$ws.Range("H106").Value = 1764.2858

$ws.Range("H112").Value = 3520
$ws.Range("J112").Value = 3483.077
$ws.Range("L112").Value = 10449.231
$ws.Range("N112").Value = -12665.231

$ws.Range("H138").Value = 2401.611
$ws.Range("I138").Value = 2248.3845
$ws.Range("K138").Value = 6745.1535
$ws.Range("M138").Value = -1605.1535

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1202.1786
$ws.Range("I32").Value = 1230.321
$ws.Range("K32").Value = 1230.321
$ws.Range("M32").Value = -943.3209999999999

$ws.Range("H45").Value = 5784.3076
$ws.Range("I45").Value = 3149.75
$ws.Range("K45").Value = 3149.75
$ws.Range("M45").Value = -2772.75

$ws.Range("H110").Value = 5542.5625
$ws.Range("I110").Value = 3532
$ws.Range("J110").Value = 8127.5713
$ws.Range("K110").Value = 3532
$ws.Range("L110").Value = 8127.5713
$ws.Range("M110").Value = -1487
$ws.Range("N110").Value = -12217.5713

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 3000
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

$ws.Range("H25").Value = 14
$ws.Range("I25").Value = 14
$ws.Range("K25").Value = 14
$ws.Range("M25").Value = 221

$ws.Range("H105").Value = 2330.5
$ws.Range("I105").Value = 2330.5
$ws.Range("K105").Value = 2330.5
$ws.Range("M105").Value = -583.5

$ws.Range("H117").Value = 38849.668
$ws.Range("J117").Value = 38849.668
$ws.Range("L117").Value = 38849.668
$ws.Range("N117").Value = -48027.668

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 41017.6
$ws.Range("J74").Value = 41017.6
$ws.Range("L74").Value = 41017.6
$ws.Range("N74").Value = -42765.6

$ws.Range("H77").Value = 41017.6
$ws.Range("J77").Value = 41017.6
$ws.Range("L77").Value = 123052.8
$ws.Range("N77").Value = -131788.8

$ws.Range("H112").Value = 79921.336
$ws.Range("J112").Value = 79921.336
$ws.Range("L112").Value = 79921.336
$ws.Range("N112").Value = -82875.336

$ws.Range("H134").Value = 6031.304
$ws.Range("I134").Value = 5260.6665
$ws.Range("K134").Value = 15781.9995
$ws.Range("M134").Value = -13246.9995

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1004.7778
$ws.Range("I12").Value = 180
$ws.Range("J12").Value = 1107.875
$ws.Range("K12").Value = 540
$ws.Range("L12").Value = 3323.625
$ws.Range("M12").Value = -367
$ws.Range("N12").Value = -3669.625

$ws.Range("H126").Value = 7565.6
$ws.Range("I126").Value = 5707.25
$ws.Range("K126").Value = 17121.75
$ws.Range("M126").Value = -12181.75

$ws.Range("H140").Value = 948829.6
$ws.Range("I140").Value = 1682.5
$ws.Range("K140").Value = 5047.5
$ws.Range("M140").Value = 132.5

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 8536.5
$ws.Range("J2").Value = 20209.8
$ws.Range("L2").Value = 20209.8
$ws.Range("N2").Value = -20435.8

$ws.Range("H59").Value = 34999.75
$ws.Range("J59").Value = 34999.75
$ws.Range("L59").Value = 34999.75
$ws.Range("N59").Value = -36165.75

$ws.Range("H99").Value = 16999.8
$ws.Range("I99").Value = 11749.75
$ws.Range("J99").Value = 38000
$ws.Range("K99").Value = 11749.75
$ws.Range("L99").Value = 38000
$ws.Range("M99").Value = -9503.75
$ws.Range("N99").Value = -42492

$ws.Range("H132").Value = 3361.8333
$ws.Range("I132").Value = 1313.2858
$ws.Range("K132").Value = 3939.8574
$ws.Range("M132").Value = -1409.8574

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 594.64703
$ws.Range("I55").Value = 629.2857
$ws.Range("K55").Value = 629.2857
$ws.Range("M55").Value = -456.2857

$ws.Range("H100").Value = 5345.769
$ws.Range("J100").Value = 6999.8335
$ws.Range("L100").Value = 6999.8335
$ws.Range("N100").Value = -8081.8335

$ws.Range("H132").Value = 9394.157999999999
$ws.Range("I132").Value = 10988.407
$ws.Range("K132").Value = 32965.221
$ws.Range("M132").Value = -30435.221

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 12000
$ws.Range("J33").Value = 12000
$ws.Range("L33").Value = 12000
$ws.Range("N33").Value = -12500

$ws.Range("H36").Value = 12000
$ws.Range("J36").Value = 12000
$ws.Range("L36").Value = 12000
$ws.Range("N36").Value = -12500

$ws.Range("H124").Value = 35000
$ws.Range("J124").Value = 35000
$ws.Range("L124").Value = 35000
$ws.Range("N124").Value = -44820

$ws.Range("H126").Value = 4293.2
$ws.Range("I126").Value = 4024
$ws.Range("K126").Value = 12072
$ws.Range("M126").Value = -9602

$ws.Range("H132").Value = 2796.0164
$ws.Range("I132").Value = 2488.691
$ws.Range("K132").Value = 7466.072999999999
$ws.Range("M132").Value = -4936.072999999999
